$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the six "JUAN EDGARDO REBOLLEDO JIMENEZ" detail rows (old rows 17-22)
$ws.Range("B17:J22").EntireRow.Delete()

# Remove the now-shifted old row 16 (1048600585 / period 2103), keeping the
# row that used to be 23 (20352225 / period 2203) as the sole detail row (16)
$ws.Range("B16:J16").EntireRow.Delete()

# Update the summary figures to match the single remaining worker/period
$ws.Range("E11").Value = 76662
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
